$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Dynamic the product, size and color" - update the Size and Color values for the
# single data row (size bumped from 28 to 29, color switched from Blue to Green).
$ws.Range("M2").Value = 29
$ws.Range("N2").Value = "Green"

# Let Excel recompute the "best fit" column widths for the sheet now that the
# data/headers have changed (mirrors what Excel does automatically on save).
for ($col = 1; $col -le 14; $col++) {
    $ws.Columns.Item($col).AutoFit() | Out-Null
}

# "Also placed the order" - move the active selection to L4.
$ws.Range("L4").Select() | Out-Null
